$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Several match rows were re-sorted (the header/id column A keeps its running
# index, but the underlying match data in columns B:AC moved between the two
# rows of each pair). Swap the B:AC payload between each pair of rows.
function Swap-RowData($ws, $r1, $r2) {
    for ($c = 2; $c -le 29; $c++) {
        $cell1 = $ws.Cells.Item($r1, $c)
        $cell2 = $ws.Cells.Item($r2, $c)
        $v1 = $cell1.Value()
        $v2 = $cell2.Value()
        $cell1.Value = $v2
        $cell2.Value = $v1
    }
}

Swap-RowData $ws 6 7
Swap-RowData $ws 16 17
Swap-RowData $ws 20 21
Swap-RowData $ws 86 87
Swap-RowData $ws 125 126

# Append the new match as row 138, copying formatting from the last existing
# data row (137) so the id cell keeps its bold/bordered style and the date
# cell keeps its custom date/time number format.
$ws.Range("A137").Copy($ws.Range("A138"))
$ws.Range("E137").Copy($ws.Range("E138"))

$ws.Cells.Item(138, 1).Value = 136
$ws.Cells.Item(138, 2).Value = 8091145
$ws.Cells.Item(138, 3).Value = "Germany Verbandsliga"
$ws.Cells.Item(138, 4).Value = "Germany Verbandsliga"
$ws.Cells.Item(138, 5).Value = 45396.41666666666
$ws.Cells.Item(138, 6).Value = "RotWeiss Darmstadt"
$ws.Cells.Item(138, 7).Value = "SV Eintracht WaldMichelbach"
$ws.Cells.Item(138, 8).Value = 1
$ws.Cells.Item(138, 9).Value = 1
$ws.Cells.Item(138, 10).Value = "D"
$ws.Cells.Item(138, 11).Value = 1.909
$ws.Cells.Item(138, 12).Value = 4
$ws.Cells.Item(138, 13).Value = 2.9
$ws.Cells.Item(138, 14).Value = 1.909
$ws.Cells.Item(138, 15).Value = 4
$ws.Cells.Item(138, 16).Value = 2.9
$ws.Cells.Item(138, 17).Value = -0.5
$ws.Cells.Item(138, 18).Value = 1.975
$ws.Cells.Item(138, 19).Value = 1.825
$ws.Cells.Item(138, 20).Value = 3.75
$ws.Cells.Item(138, 21).Value = 2
$ws.Cells.Item(138, 22).Value = 1.8
$ws.Cells.Item(138, 23).Value = -1
$ws.Cells.Item(138, 24).Value = 3
$ws.Cells.Item(138, 25).Value = -1
$ws.Cells.Item(138, 26).Value = -1
$ws.Cells.Item(138, 27).Value = 0.825
$ws.Cells.Item(138, 28).Value = -1
$ws.Cells.Item(138, 29).Value = 0.8
